$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D5").Value = 6
$ws.Range("D6").Value = 1

$ws.Range("D6").Select()
